$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-29 20:59:35"

for ($row = 2; $row -le 37; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
